$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031148601803523
$ws.Range("D2").Value = 1.040909313596436
$ws.Range("E2").Value = 1.030764055098374
$ws.Range("F2").Value = 1.049117848103666
$ws.Range("I2").Value = 1.035216852450105
$ws.Range("J2").Value = 1.036285924666486
$ws.Range("K2").Value = 1.043690326243687
$ws.Range("L2").Value = 1.033574122881548
$ws.Range("M2").Value = 1.051875798331636
$ws.Range("N2").Value = 1.016088676407724

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032060128394431
$ws.Range("D3").Value = 1.041767317166288
$ws.Range("E3").Value = 1.031536809732144
$ws.Range("F3").Value = 1.050063457818646
$ws.Range("I3").Value = 1.035365694545445
$ws.Range("J3").Value = 1.036839397791536
$ws.Range("K3").Value = 1.044358724687893
$ws.Range("L3").Value = 1.034155396247195
$ws.Range("M3").Value = 1.052633245090698
$ws.Range("N3").Value = 1.016272192556716

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032650115229043
$ws.Range("D4").Value = 1.042322927955699
$ws.Range("E4").Value = 1.032037373449064
$ws.Range("F4").Value = 1.05067576123919
$ws.Range("I4").Value = 1.035460075148996
$ws.Range("J4").Value = 1.037197078348086
$ws.Range("K4").Value = 1.044791007081779
$ws.Range("L4").Value = 1.034531400701405
$ws.Range("M4").Value = 1.053123170322535
$ws.Range("N4").Value = 1.016390765544594

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032898184855489
$ws.Range("D5").Value = 1.042556607075932
$ws.Range("E5").Value = 1.03224793822235
$ws.Range("F5").Value = 1.050933275260189
$ws.Range("I5").Value = 1.035499289895941
$ws.Range("J5").Value = 1.037347337617059
$ws.Range("K5").Value = 1.044972685666293
$ws.Range("L5").Value = 1.034689443543587
$ws.Range("M5").Value = 1.05332908725958
$ws.Range("N5").Value = 1.016440571537063

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032939839101719
$ws.Range("D6").Value = 1.042595848671116
$ws.Range("E6").Value = 1.032283300425774
$ws.Range("F6").Value = 1.0509765189168
$ws.Range("I6").Value = 1.03550584705979
$ws.Range("J6").Value = 1.037372560349635
$ws.Range("K6").Value = 1.045003187160994
$ws.Range("L6").Value = 1.034715977873944
$ws.Range("M6").Value = 1.053363658788038
$ws.Range("N6").Value = 1.016448931702997

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032653429795038
$ws.Range("D7").Value = 1.042326049995056
$ws.Range("E7").Value = 1.032040186525576
$ws.Range("F7").Value = 1.050679201756504
$ws.Range("I7").Value = 1.035460600958056
$ws.Range("J7").Value = 1.03719908655192
$ws.Range("K7").Value = 1.044793434889427
$ws.Range("L7").Value = 1.034533512594645
$ws.Range("M7").Value = 1.05312592198354
$ws.Range("N7").Value = 1.016391431220574

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03145662099909
$ws.Range("D8").Value = 1.04119919104253
$ws.Range("E8").Value = 1.03102509851141
$ws.Range("F8").Value = 1.049437331284884
$ws.Range("I8").Value = 1.035267553530297
$ws.Range("J8").Value = 1.036473066782729
$ws.Range("K8").Value = 1.043916258650825
$ws.Range("L8").Value = 1.033770590760079
$ws.Range("M8").Value = 1.052131820281672
$ws.Range("N8").Value = 1.016150732397332

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029349027333419
$ws.Range("D9").Value = 1.039216840087114
$ws.Range("E9").Value = 1.029240581309583
$ws.Range("F9").Value = 1.047252359113965
$ws.Range("I9").Value = 1.034912629995352
$ws.Range("J9").Value = 1.035190299590403
$ws.Range("K9").Value = 1.042368959704436
$ws.Range("L9").Value = 1.032425368001598
$ws.Range("M9").Value = 1.050378661939392
$ws.Range("N9").Value = 1.015725272591836

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027944926820261
$ws.Range("D10").Value = 1.03789758722709
$ws.Range("E10").Value = 1.028053809912794
$ws.Range("F10").Value = 1.045798055865456
$ws.Range("I10").Value = 1.034666142231606
$ws.Range("J10").Value = 1.034332878134033
$ws.Range("K10").Value = 1.041336419618208
$ws.Range("L10").Value = 1.0315280416173
$ws.Range("M10").Value = 1.049209008564177
$ws.Range("N10").Value = 1.015440768915799

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027337176800189
$ws.Range("D11").Value = 1.037326902308375
$ws.Range("E11").Value = 1.027540631781856
$ws.Range("F11").Value = 1.045168899985498
$ws.Range("I11").Value = 1.034557077912625
$ws.Range("J11").Value = 1.033961085662616
$ws.Range("K11").Value = 1.040889093776141
$ws.Range("L11").Value = 1.031139381552834
$ws.Range("M11").Value = 1.048702342428471
$ws.Range("N11").Value = 1.015317375367554

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027111467626357
$ws.Range("D12").Value = 1.037115010040097
$ws.Range("E12").Value = 1.027350121388362
$ws.Range("F12").Value = 1.044935290013867
$ws.Range("I12").Value = 1.034516216676168
$ws.Range("J12").Value = 1.03382290773254
$ws.Range("K12").Value = 1.040722903725427
$ws.Range("L12").Value = 1.030995000148485
$ws.Range("M12").Value = 1.048514115457716
$ws.Range("N12").Value = 1.015271511579262

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027159881354631
$ws.Range("D13").Value = 1.037160457759407
$ws.Range("E13").Value = 1.027390981658687
$ws.Range("F13").Value = 1.044985396208311
$ws.Range("I13").Value = 1.03452499737246
$ws.Range("J13").Value = 1.033852550865972
$ws.Range("K13").Value = 1.040758553545212
$ws.Range("L13").Value = 1.031025971145828
$ws.Range("M13").Value = 1.048554492057045
$ws.Range("N13").Value = 1.015281350866275

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02731851886449
$ws.Range("D14").Value = 1.03730938547503
$ws.Range("E14").Value = 1.027524881946323
$ws.Range("F14").Value = 1.045149587939303
$ws.Range("I14").Value = 1.034553707443807
$ws.Range("J14").Value = 1.03394966541029
$ws.Range("K14").Value = 1.040875357127908
$ws.Range("L14").Value = 1.031127447256187
$ws.Range("M14").Value = 1.048686784104952
$ws.Range("N14").Value = 1.015313584861032

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027416265501721
$ws.Range("D15").Value = 1.037401156128283
$ws.Range("E15").Value = 1.02760739652809
$ws.Range("F15").Value = 1.045250763390751
$ws.Range("I15").Value = 1.034571350322695
$ws.Range("J15").Value = 1.034009490625546
$ws.Range("K15").Value = 1.04094731928881
$ws.Range("L15").Value = 1.031189967980989
$ws.Range("M15").Value = 1.048768289847182
$ws.Range("N15").Value = 1.015333441333573

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027985266216847
$ws.Range("D16").Value = 1.037935473624937
$ws.Range("E16").Value = 1.028087882758237
$ws.Range("F16").Value = 1.045839822930705
$ws.Range("I16").Value = 1.03467333136616
$ws.Range("J16").Value = 1.034357541833081
$ws.Range("K16").Value = 1.041366102427985
$ws.Range("L16").Value = 1.031553833439437
$ws.Range("M16").Value = 1.049242630248439
$ws.Range("N16").Value = 1.015448953922475

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028342248873741
$ws.Range("D17").Value = 1.038270787826999
$ws.Range("E17").Value = 1.028389467927394
$ws.Range("F17").Value = 1.04620947717773
$ws.Range("I17").Value = 1.034736677126216
$ws.Range("J17").Value = 1.034575725738841
$ws.Range("K17").Value = 1.041628733459818
$ws.Range("L17").Value = 1.03178204743267
$ws.Range("M17").Value = 1.049540119163671
$ws.Range("N17").Value = 1.015521358200099

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02855049340442
$ws.Range("D18").Value = 1.038466424896524
$ws.Range("E18").Value = 1.028565445066679
$ws.Range("F18").Value = 1.046425144865495
$ws.Range("I18").Value = 1.034773400442445
$ws.Range("J18").Value = 1.034702938161763
$ws.Range("K18").Value = 1.041781899511166
$ws.Range("L18").Value = 1.031915149934933
$ws.Range("M18").Value = 1.049713620180722
$ws.Range("N18").Value = 1.015563570917108

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028621503209322
$ws.Range("D19").Value = 1.038533141194514
$ws.Range("E19").Value = 1.028625460149444
$ws.Range("F19").Value = 1.046498691197055
$ws.Range("I19").Value = 1.034785883926457
$ws.Range("J19").Value = 1.034746305705401
$ws.Range("K19").Value = 1.041834121379375
$ws.Range("L19").Value = 1.031960532574714
$ws.Range("M19").Value = 1.049772776244408
$ws.Range("N19").Value = 1.015577961048938

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028303945669784
$ws.Range("D20").Value = 1.038234806194518
$ws.Range("E20").Value = 1.028357103689369
$ws.Range("F20").Value = 1.04616981111606
$ws.Range("I20").Value = 1.034729904013515
$ws.Range("J20").Value = 1.034552321893011
$ws.Range("K20").Value = 1.04160055792234
$ws.Range("L20").Value = 1.031757563346572
$ws.Range("M20").Value = 1.049508203399119
$ws.Range("N20").Value = 1.015513591917961

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027271803063466
$ws.Range("D21").Value = 1.037265527603147
$ws.Range("E21").Value = 1.027485448689286
$ws.Range("F21").Value = 1.045101235167307
$ws.Range("I21").Value = 1.034545262696612
$ws.Range("J21").Value = 1.033921069733306
$ws.Range("K21").Value = 1.040840962294808
$ws.Range("L21").Value = 1.031097565497329
$ws.Range("M21").Value = 1.048647828174144
$ws.Range("N21").Value = 1.015304093576372

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.026623064216575
$ws.Range("D22").Value = 1.036656599138668
$ws.Range("E22").Value = 1.026938023269053
$ws.Range("F22").Value = 1.044429880386907
$ws.Range("I22").Value = 1.034427147340552
$ws.Range("J22").Value = 1.033523727652655
$ws.Range("K22").Value = 1.040363182169199
$ws.Range("L22").Value = 1.030682508037817
$ws.Range("M22").Value = 1.04810671119908
$ws.Range("N22").Value = 1.015172200660536

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026966952707311
$ws.Range("D23").Value = 1.036979356227516
$ws.Range("E23").Value = 1.027228164736436
$ws.Range("F23").Value = 1.044785730289708
$ws.Range("I23").Value = 1.034489954151901
$ws.Range("J23").Value = 1.033734408424825
$ws.Range("K23").Value = 1.040616480301632
$ws.Range("L23").Value = 1.030902546090776
$ws.Range("M23").Value = 1.048393582774315
$ws.Range("N23").Value = 1.015242135870514

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028321253169169
$ws.Range("D24").Value = 1.038251064578788
$ws.Range("E24").Value = 1.028371727486656
$ws.Range("F24").Value = 1.046187734333262
$ws.Range("I24").Value = 1.034732965187461
$ws.Range("J24").Value = 1.034562897239825
$ws.Range("K24").Value = 1.041613289303565
$ws.Range("L24").Value = 1.031768626684819
$ws.Range("M24").Value = 1.049522624818844
$ws.Range("N24").Value = 1.015517101226762

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029893725904318
$ws.Range("D25").Value = 1.039728923466168
$ws.Range("E25").Value = 1.029701415616014
$ws.Range("F25").Value = 1.047816819645062
$ws.Range("I25").Value = 1.035006129496208
$ws.Range("J25").Value = 1.035522325521636
$ws.Range("K25").Value = 1.042769156237689
$ws.Range("L25").Value = 1.032773235174733
$ws.Range("M25").Value = 1.050832056271856
$ws.Range("N25").Value = 1.015835418017196
